$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.033.54'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '1.828.48'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('D4').Value = "'0.9990"
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '240.71'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '0.6211'
$ws.Range('E6').Value = '  -6.22%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '44.44'
$ws.Range('E8').Value = '  +6.12%  '
$ws.Range('D9').Value = '0.07457'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').Value = '0.2917'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').Value = '22.65'
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '0.07604'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').Value = '1.828.57'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '4.954'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '0.6621'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '81.98'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').Value = "'0.000009174"
$ws.Range('E17').Value = '  +10.23%  '
$ws.Range('D18').Value = '6.001'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').Value = '29.038.44'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '225.01'
$ws.Range('E20').Value = '  -0.61%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '12.33'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'1.000"
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '7.167'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = "'1.000"
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '159.38'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '8.396'
$ws.Range('E26').Value = '  -2.19%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '0.1352'
$ws.Range('E27').Value = '  -3.13%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '17.79'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '1.498'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = "'4.050"
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '4.026'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '1.202'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.05229'
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '1.835'
$ws.Range('E34').Value = '  -1.56%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.7344'
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.148'
$ws.Range('E36').Value = '  +1.28%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '2.647'
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.276.09'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = "'2.750"
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.01780"
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '6.313'
$ws.Range('E41').Value = '  +7.05%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8936'
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = "'1.000"
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '101.71'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.977.99'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = '0.5111'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '63.44'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = "'0.00000000119"
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.704'
$ws.Range('E49').Value = '  -3.05%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = '0.3959'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '8.808'
$ws.Range('E51').Value = '  -0.28%  '
